# Grafikplanering.xlsx update:
#  - Row 89 ("6.3 Menybakgrund" / "0% Färdigt") is removed entirely, shifting
#    every row below it up by one (old row 90 -> 89, ..., old row 118 -> 117).
#  - The completion percentage for "6.2 Huvudmeny" (row 87, column I) is
#    updated from "30% Färdigt" to "70% Färdigt".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "6.3 Menybakgrund" row; everything below shifts up one row.
$ws.Rows("89").Delete()

# Bump the "6.2 Huvudmeny" progress value.
$ws.Range("I87").Value = "70% Färdigt"

# Restore the selection the author left the sheet on.
$ws.Range("S93").Select() | Out-Null
